# Re-pull data / push all data / mean calculation
# Updates the "dSF" (column F) values that changed on the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    5  = -11
    6  = -2
    7  = 0
    8  = -3
    15 = -4
    17 = -3
    18 = 0
    26 = 1
    27 = -1
    29 = -1
    35 = 1
    36 = -4
    37 = 7
    39 = -2
    45 = -2
    49 = -5
    51 = -2
    52 = 3
    56 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
